# Update the build timestamp embedded in the version string across the workbook.
$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Jude Coal Mine, China, M2011, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 7; $row++) {
    $cell = $data.Cells.Item($row, 19)  # column S = 19
    if ($cell.Text -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
